# Applies the "Add files via upload" commit:
#  - Nunavut: appends 4 new data rows (21-24), becomes the active/selected
#    sheet with active cell H9.
#  - Prince Edward Island: appends 5 new data rows (14-18), active cell
#    becomes L21 (no longer the active sheet).
#  - New Brunswick: loses tabSelected (no longer active sheet), active
#    cell stays C18.

$wb = $excel.ActiveWorkbook

# --- Prince Edward Island: add rows 14-18 -------------------------------
$wsPEI = $wb.Worksheets.Item("Prince Edward Island")
$wsPEI.Activate()

$peiRows = @(
    @(44548, 1, 3, 1),
    @(44550, 3, 3, 1),
    @(44579, 3, 3, 2),
    @(44592, 1, 3, 1),
    @(44609, 1, 2, 1)
)

$destRow = 14
foreach ($r in $peiRows) {
    # Copy the formatting of the last existing data row (row 13) down onto
    # the new row so the new cells pick up the same style indices.
    $wsPEI.Range("A13:D13").Copy()
    $wsPEI.Range("A$destRow`:D$destRow").PasteSpecial(-4122)

    $wsPEI.Cells.Item($destRow, 1).Value = $r[0]
    $wsPEI.Cells.Item($destRow, 2).Value = $r[1]
    $wsPEI.Cells.Item($destRow, 3).Value = $r[2]
    $wsPEI.Cells.Item($destRow, 4).Value = $r[3]

    $destRow = $destRow + 1
}

$wsPEI.Range("L21").Select()

# --- Nunavut: add rows 21-24, then make it the selected/active sheet ---
$wsNun = $wb.Worksheets.Item("Nunavut")
$wsNun.Activate()

$nunRows = @(
    @(44550, 3, 2, 1),
    @(44554, 3, 3, 3),
    @(44574, 3, 2, 1),
    @(44585, 1, 2, 1)
)

$destRow = 21
foreach ($r in $nunRows) {
    # Copy the formatting of the last existing data row (row 20) down onto
    # the new row so the new cells pick up the same style indices.
    $wsNun.Range("A20:D20").Copy()
    $wsNun.Range("A$destRow`:D$destRow").PasteSpecial(-4122)

    $wsNun.Cells.Item($destRow, 1).Value = $r[0]
    $wsNun.Cells.Item($destRow, 2).Value = $r[1]
    $wsNun.Cells.Item($destRow, 3).Value = $r[2]
    $wsNun.Cells.Item($destRow, 4).Value = $r[3]

    $destRow = $destRow + 1
}

$wsNun.Range("H9").Select()

# Nunavut ends up as the active/selected tab (matches the workbook's new
# activeTab), which also clears tabSelected from the previously-active
# New Brunswick sheet.
